$d = $word.ActiveDocument
$paragraphs = $d.Paragraphs

# Paragraphs 50-55 (1-indexed) make up the "KEY ACHIEVEMENTS AND IMPACT" bullet list.
# 1) Rewrite the first four bullets in place (text replaced, same paragraph objects).
#    Note: no trailing `r here - Range.Text already ends right before the existing
#    paragraph mark, so appending `r would insert an extra empty paragraph.
$paragraphs.Item(50).Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
$paragraphs.Item(51).Range.Text = "• Real-time collaboration at national scale"
$paragraphs.Item(52).Range.Text = "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
$paragraphs.Item(53).Range.Text = "• 23% conversion rate improvement"

# 2) Remove the last two bullet paragraphs entirely (paragraph indices are unchanged
#    because step 1 only replaced text, not paragraph marks).
$start = $paragraphs.Item(54).Range.Start
$end = $paragraphs.Item(55).Range.End
$d.Range($start, $end).Delete()

$d.Save()
